$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 86, shifting existing rows 86-146 down to 87-147
$ws.Rows.Item(86).Insert()

# Populate the new row 86 with the new record (matches the rest of the table's shape)
$ws.Range("A86").Value = 10
$ws.Range("B86").Value = "Vega Modelo de Temuco"
$ws.Range("C86").Value = "La Araucanía"
$ws.Range("D86").Value = 44566
$ws.Range("E86").Value = 9
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100103
$ws.Range("H86").Value = "Frutos de hueso (carozo)"
$ws.Range("I86").Value = 100103002
$ws.Range("J86").Value = "Ciruela"
$ws.Range("K86").Value = "Blackbeaut"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 125
$ws.Range("N86").Value = 15000
$ws.Range("O86").Value = 15000
$ws.Range("P86").Value = 15000
$ws.Range("Q86").Value = "$/caja 18 kilos granel"
$ws.Range("R86").Value = "Región de O'Higgins"
$ws.Range("S86").Value = 833
$ws.Range("T86").Value = 18

# Apply the same date style (numFmtId 165) used by the rest of column D
$ws.Range("D86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
